$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 54
$ws.Range("C3").Value = 55
$ws.Range("B4").Value = "<shift>"
$ws.Range("B5").Value = "<of>"
$ws.Range("C5").Value = 55
$ws.Range("C6").Value = 58
$ws.Range("B7").Value = "<which>"
$ws.Range("C7").Value = 54
$ws.Range("C8").Value = 47
$ws.Range("B9").Value = "<delta>"
$ws.Range("C9").Value = 15
